# Add 5 topic textboxes ("@DEV", "@INTEG", "@TEST", "@RUN", "@PERF") to the
# workshop-overview slide (slide 2), matching the shape ids 51-55 used by
# the canonical OOXML.
#
# This runtime assigns a new shape's Id by filling the smallest gap missing
# from the ids that have ever been issued on the slide in this session (a
# deleted shape's id is not reclaimed/reused). Slide 2 already uses ids
# 1,3,4,5,6,9,10,12,19,20,22,23,25,26,27,28,30,31,33,35,36,37,38,39,40,42,
# 43,46,47,48,49,50 -- so the natural gap-filling sequence for new shapes
# is 2,7,8,11,13,14,15,16,17,18,21,24,29,32,34,41,44,45,51,52,53,... . To
# land the 5 real textboxes on ids 51-55 (as in the target deck) we first
# "burn" through the first 18 gap ids with disposable textboxes that are
# immediately deleted.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

for ($i = 0; $i -lt 18; $i++) {
    $burn = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
    $burn.Delete()
}

function Add-TopicTextbox($leftEmu, $topEmu, $widthEmu, $heightEmu, $text) {
    $left = $leftEmu / 12700
    $top = $topEmu / 12700
    $width = $widthEmu / 12700
    $height = $heightEmu / 12700

    $shp = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
    $shp.TextFrame.WordWrap = $false
    $tr = $shp.TextFrame.TextRange
    $tr.Text = $text
    $tr.Font.Bold = $true
    $tr.Font.Name = "Calibri"
    $tr.Font.NameComplexScript = "Calibri"
    $shp.TextFrame.AutoSize = 1
    $shp.Fill.Visible = $false
    $shp
}

Add-TopicTextbox 2484251 711131  787395 369332 "@DEV"   | Out-Null
Add-TopicTextbox 4947689 711131  979755 369332 "@INTEG" | Out-Null
Add-TopicTextbox 7402100 711131  842335 369332 "@TEST"  | Out-Null
Add-TopicTextbox 2493278 3603255 824753 369332 "@RUN"   | Out-Null
Add-TopicTextbox 4947689 3603255 863412 369332 "@PERF"  | Out-Null
